# Update the Kng1-Plaur TPM-derived metrics for rows 2-6 with new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 14.349718
$ws.Range("N2").Value = 43.049154
$ws.Range("O2").Value = 0.1016415840981481
$ws.Range("P2").Value = 0.1034081666702025
$ws.Range("Q2").Value = 2.176144301178667
$ws.Range("R2").Value = 19.585298710608
$ws.Range("S2").Value = 0.1016415840981481
$ws.Range("T2").Value = 0.1034081666702025

# Row 3
$ws.Range("O3").Value = 0.04778708884009916
$ws.Range("P3").Value = 0.04861765281706964
$ws.Range("S3").Value = 0.04778708884009916
$ws.Range("T3").Value = 0.04861765281706964

# Row 4
$ws.Range("M4").Value = 66.43651233333334
$ws.Range("N4").Value = 199.309537
$ws.Range("O4").Value = 0.4705815372480596
$ws.Range("P4").Value = 0.4787604843769264
$ws.Range("Q4").Value = 10.07514138635822
$ws.Range("R4").Value = 90.676272477224
$ws.Range("S4").Value = 0.4705815372480596
$ws.Range("T4").Value = 0.4787604843769264

# Row 5
$ws.Range("M5").Value = 7.2355625
$ws.Range("N5").Value = 14.471125
$ws.Range("O5").Value = 0.05125076564857627
$ws.Range("P5").Value = 0.03476102006337534
$ws.Range("Q5").Value = 1.097277876833334
$ws.Range("R5").Value = 6.583667261
$ws.Range("S5").Value = 0.05125076564857627
$ws.Range("T5").Value = 0.03476102006337534

# Row 6
$ws.Range("M6").Value = 46.41124333333334
$ws.Range("N6").Value = 139.23373
$ws.Range("O6").Value = 0.328739024165117
$ws.Range("P6").Value = 0.3344526760724259
$ws.Range("Q6").Value = 7.038295992328891
$ws.Range("R6").Value = 63.34466393096001
$ws.Range("S6").Value = 0.328739024165117
$ws.Range("T6").Value = 0.3344526760724259
